# Updated cryptos list - apply cell-level changes as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$ref, [string]$val) {
    $ws.Range($ref).Value = "'" + $val
    $ws.Range($ref).Style = "Normal"
}

Set-TextCell 'D2' '30.282.45'
Set-TextCell 'E2' '  +2.01%  '
Set-TextCell 'D3' '2.089.56'
Set-TextCell 'E3' '  -0.38%  '
Set-TextCell 'D4' '1.003'
Set-TextCell 'E4' '  -0.31%  '
Set-TextCell 'D5' '343.26'
Set-TextCell 'E5' '  -0.22%  '
Set-TextCell 'D6' '1.003'
Set-TextCell 'E6' '  -0.25%  '
Set-TextCell 'D7' '0.5211'
Set-TextCell 'E7' '  +1.43%  '
Set-TextCell 'D8' '0.4402'
Set-TextCell 'E8' '  -0.15%  '
Set-TextCell 'D9' '54.47'
Set-TextCell 'E9' '  +2.99%  '
Set-TextCell 'D10' '0.09333'
Set-TextCell 'E10' '  +1.90%  '
Set-TextCell 'D11' '1.169'
Set-TextCell 'E11' '  -0.31%  '
Set-TextCell 'D12' '24.77'
Set-TextCell 'E12' '  -0.87%  '
Set-TextCell 'D13' '8.644'
Set-TextCell 'E13' '  +4.67%  '
Set-TextCell 'D14' '6.897'
Set-TextCell 'E14' '  +2.01%  '
Set-TextCell 'D15' '2.089.63'
Set-TextCell 'E15' '  -0.93%  '
Set-TextCell 'D16' '100.98'
Set-TextCell 'E16' '  +1.34%  '
Set-TextCell 'D17' '0.00001156'
Set-TextCell 'E17' '  +0.55%  '
Set-TextCell 'D18' '1.004'
Set-TextCell 'E18' '  -0.35%  '
Set-TextCell 'D19' '21.14'
Set-TextCell 'E19' '  +1.68%  '
Set-TextCell 'D20' '0.06677'
Set-TextCell 'E20' '  +0.77%  '
Set-TextCell 'D21' '6.351'
Set-TextCell 'E21' '  +2.63%  '
Set-TextCell 'D22' '1.003'
Set-TextCell 'E22' '  -0.27%  '
Set-TextCell 'D23' '30.265.14'
Set-TextCell 'E23' '  +1.73%  '
Set-TextCell 'D24' '12.52'
Set-TextCell 'E24' '  -0.56%  '
Set-TextCell 'D25' '2.293'
Set-TextCell 'E25' '  -0.98%  '
Set-TextCell 'D26' '21.77'
Set-TextCell 'E26' '  -0.58%  '
Set-TextCell 'D27' '162.24'
Set-TextCell 'E27' '  +0.11%  '
Set-TextCell 'D28' '2.513'
Set-TextCell 'E28' '  -0.61%  '
Set-TextCell 'D29' '132.84'
Set-TextCell 'E29' '  +0.01%  '
Set-TextCell 'D30' '1.128'
Set-TextCell 'E30' '  -0.18%  '
Set-TextCell 'B31' 'Stellar'
Set-TextCell 'C31' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D31' '0.1047'
Set-TextCell 'E31' '  -0.16%  '
Set-TextCell 'B32' 'ARBITRUM'
Set-TextCell 'C32' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D32' '1.658'
Set-TextCell 'E32' '  -0.27%  '
Set-TextCell 'D33' '6.212'
Set-TextCell 'E33' '  +0.72%  '
Set-TextCell 'D34' '6.636'
Set-TextCell 'E34' '  +10.15%  '
Set-TextCell 'D35' '3.857'
Set-TextCell 'E35' '  -2.15%  '
Set-TextCell 'D36' '10.16'
Set-TextCell 'E36' '  -2.70%  '
Set-TextCell 'D37' '0.02628'
Set-TextCell 'E37' '  +2.48%  '
Set-TextCell 'D38' '0.06797'
Set-TextCell 'E38' '  +1.11%  '
Set-TextCell 'D39' '0.6978'
Set-TextCell 'E39' '  +1.66%  '
Set-TextCell 'E40' '  +4.04%  '
Set-TextCell 'D41' '12.51'
Set-TextCell 'E41' '  +0.54%  '
Set-TextCell 'D42' '0.2209'
Set-TextCell 'E42' '  -1.00%  '
Set-TextCell 'D43' '0.6803'
Set-TextCell 'E43' '  +2.12%  '
Set-TextCell 'D44' '14.36'
Set-TextCell 'E44' '  +0.83%  '
Set-TextCell 'D45' '2.333'
Set-TextCell 'E45' '  +1.56%  '
Set-TextCell 'E46' '  -0.19%  '
Set-TextCell 'D47' '1.371'
Set-TextCell 'E47' '  +18.25%  '
Set-TextCell 'D48' '3.631'
Set-TextCell 'E49' '  -0.65%  '
Set-TextCell 'D50' '1.216'
Set-TextCell 'E50' '  +8.71%  '
Set-TextCell 'D51' '1.215'
Set-TextCell 'E51' '  -0.44%  '
